$d = $word.ActiveDocument

$replacements = @(
    @("51×78=", "81×74="),
    @("72×88=", "27×11="),
    @("88×42=", "32×38="),
    @("73×40=", "51×94="),
    @("48×60=", "75×66="),
    @("11×95=", "81×91="),
    @("33×24=", "29×54="),
    @("59×56=", "71×62="),
    @("35×84=", "39×26="),
    @("41×16=", "31×63="),
    @("63×72=", "73×28="),
    @("35×97=", "48×42="),
    @("88×32=", "63×39="),
    @("47×97=", "67×65="),
    @("60×34=", "14×62="),
    @("84×47=", "88×80="),
    @("44×16=", "67×92="),
    @("86×30=", "54×51="),
    @("69×12=", "35×63="),
    @("53×18=", "64×24="),
    @("81×96=", "23×81="),
    @("25×76=", "34×31="),
    @("65×58=", "27×57="),
    @("85×76=", "18×39="),
    @("40×82=", "57×17=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
